# The source data (weekly Cilantro price report) gained one more week's
# worth of observations. A new record was inserted as row 26, pushing the
# previously-existing rows 26-59 down to rows 27-60 (their values are
# unchanged - only their row position moved).
#
# The new row 26 carries its own data: 2021-12-27 (serial 44557), volume
# 300, min/max/avg prices 2800/3000/2900, and $/Kg price 1450 - the same
# figures that (coincidentally) already appear on the last row of the
# table, consistent with the most recent observation being duplicated
# going into the following week's report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; Excel shifts rows 26:59 down to
# 27:60 automatically, preserving all of their existing content.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new observation.
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44557
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112040
$ws.Cells.Item(26, 7).Value = "Cilantro"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 2800
$ws.Cells.Item(26, 12).Value = 3000
$ws.Cells.Item(26, 13).Value = 2900
$ws.Cells.Item(26, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 1450
$ws.Cells.Item(26, 17).Value = 2
$ws.Cells.Item(26, 18).Value = "Hortaliza"
